$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 70 (shifts existing rows 70.. down by one)
$ws.Rows.Item(70).Insert()

# Populate the newly inserted row 70 with the new record
$ws.Cells.Item(70, 1).Value = 10
$ws.Cells.Item(70, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(70, 3).Value = "La Araucanía"
$ws.Cells.Item(70, 4).Value = 44469
$ws.Cells.Item(70, 5).Value = 9
$ws.Cells.Item(70, 6).Value = 100112001
$ws.Cells.Item(70, 7).Value = "Berenjena"
$ws.Cells.Item(70, 8).Value = "Sin especificar"
$ws.Cells.Item(70, 9).Value = "Primera"
$ws.Cells.Item(70, 10).Value = 100
$ws.Cells.Item(70, 11).Value = 12000
$ws.Cells.Item(70, 12).Value = 12000
$ws.Cells.Item(70, 13).Value = 12000
$ws.Cells.Item(70, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(70, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(70, 16).Value = 200
$ws.Cells.Item(70, 17).Value = 60
$ws.Cells.Item(70, 18).Value = "Hortaliza"
